$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Notes")

# --- Row 3 / column C3 formatting tweak (wrap + vertical center, no horizontal) ---
$ws.Range("C3").VerticalAlignment = -4108   # xlCenter
$ws.Range("C3").WrapText = $true

# --- New reference string (C4) ---
$refText = @'
1. https://www.youtube.com/watch?v=s1fgWfyfKVQ
'@
$refText = $refText.TrimEnd("`r","`n")

# --- New rich "Notes" text (B4), two runs: bold title + regular body ---
$boldRun = 'Center of Pressure and Aerodynamic center '
$bodyRun = @'
- Center of pressure is where Lift and Drag happens, but the issue is that CP changes any time Angle of attack changes. The aerodynamic force, CP, and thus moment changes everytime the angle of attack changes. So Engineers use Aerodynamic center for convenience. Aerodynamic center (AC) is fixed, and we can apply the Lift and Drag force on the Aerodynamic center (AC). This makes it easier to analyze pitch moments. But we have to remember to add M_ac or the moment applied at the aerodynamic center ( this moment compensates us moving the Lift and drag from the CP to AC). Watch the video in Ref 1 for full explanation. Remeber the Force on a beam analogy. For moving the Force to another point, we add an additional moment on the new point.  The quantity of the moment is the distance between the original point and the new point times the force. 
For the rocket, since the change in angle of attack is not too great, which means the CP won't move as much, I can stick to finding the CP when alpha is 0 using barrowman's method and apply lift and drag at the found CP. 
'@
$bodyRun = $bodyRun.TrimEnd("`r","`n")

$noteText = $boldRun + $bodyRun

# Row 4 values
$ws.Range("B4").Value = $noteText
$ws.Range("B4").Characters(1, $boldRun.Length).Font.Bold = $true

$ws.Range("C4").Value = $refText

# Row 4 formatting: B4 center/center/wrap (matches column B default style),
# C4 left/center/wrap
$ws.Range("B4").HorizontalAlignment = -4108  # xlCenter
$ws.Range("B4").VerticalAlignment = -4108    # xlCenter
$ws.Range("B4").WrapText = $true

$ws.Range("C4").HorizontalAlignment = -4131  # xlLeft
$ws.Range("C4").VerticalAlignment = -4108    # xlCenter
$ws.Range("C4").WrapText = $true

# Row heights
$ws.Rows(3).RowHeight = 30
$ws.Rows(4).RowHeight = 210

# --- Header cell B1: bold, centered, wrap (column-wide default for B going forward) ---
$ws.Range("B1").HorizontalAlignment = -4108
$ws.Range("B1").VerticalAlignment = -4108
$ws.Range("B1").WrapText = $true
$ws.Range("B1").Font.Bold = $true

# --- Selection / view state ---
$ws.Range("B4").Select()
